# Auto-generated edit script applying cryptos.xlsx update
# Updates price (D) and volume-change (E) columns for rows 2-47,
# and replaces rows 48-51 coin listing (EnergySwap dropped, BabyDogeCoin added,
# with Maker/FLOKI/Monero shifting position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. '528.92') need to be
# forced to Text format first, otherwise Excel auto-converts them to floating
# point numbers and we lose the exact original text (trailing zeros, etc.).
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D12",
    "D14",
    "D16",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D33",
    "D34",
    "D36",
    "D37",
    "D39",
    "D40",
    "D43",
    "D45",
    "D49",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @(
    @("D2", "69.232.96"),
    @("E2", "  +1.33%  "),
    @("D3", "3.891.13"),
    @("E3", "  -0.33%  "),
    @("E4", "  -0.01%  "),
    @("D5", "528.92"),
    @("E5", "  +9.00%  "),
    @("D6", "144.40"),
    @("E6", "  -1.01%  "),
    @("E7", "  -1.64%  "),
    @("E8", "  +0.09%  "),
    @("D9", "0.718"),
    @("E9", "  -3.13%  "),
    @("E10", "  -5.77%  "),
    @("E11", "  -5.54%  "),
    @("D12", "41.96"),
    @("E12", "  -2.33%  "),
    @("D13", "4.511.54"),
    @("E13", "  -0.10%  "),
    @("D14", "10.21"),
    @("E14", "  -2.34%  "),
    @("D15", "3.911.36"),
    @("E15", "  -0.21%  "),
    @("D16", "13.99"),
    @("E16", "  -1.66%  "),
    @("E17", "  +6.68%  "),
    @("E18", "  -1.35%  "),
    @("D19", "20.01"),
    @("E19", "  +0.09%  "),
    @("D20", "69.200.02"),
    @("E20", "  +1.13%  "),
    @("D21", "423.27"),
    @("E21", "  -1.82%  "),
    @("D22", "3.40"),
    @("E22", "  -4.93%  "),
    @("D23", "14.16"),
    @("E23", "  -4.06%  "),
    @("D24", "87.55"),
    @("E24", "  -2.00%  "),
    @("D25", "4.01"),
    @("E25", "  +8.27%  "),
    @("D26", "11.37"),
    @("E26", "  -8.03%  "),
    @("E27", "  -3.69%  "),
    @("D28", "36.34"),
    @("E28", "  -2.62%  "),
    @("D29", "695.61"),
    @("E29", "  -3.02%  "),
    @("D30", "13.19"),
    @("E30", "  -1.50%  "),
    @("E31", "  -3.30%  "),
    @("E32", "  -2.43%  "),
    @("D33", "67.97"),
    @("E33", "  +10.28%  "),
    @("D34", "0.432"),
    @("E34", "  +6.91%  "),
    @("E35", "  -4.49%  "),
    @("D36", "5.91"),
    @("E36", "  -2.63%  "),
    @("D37", "40.00"),
    @("E37", "  -1.89%  "),
    @("E38", "  -0.08%  "),
    @("D39", "0.148"),
    @("E39", "  +0.42%  "),
    @("D40", "0.999"),
    @("E40", "  -0.19%  "),
    @("E41", "  +6.41%  "),
    @("E42", "  +7.14%  "),
    @("D43", "0.0481"),
    @("E43", "  -2.83%  "),
    @("E44", "  -6.96%  "),
    @("D45", "3.41"),
    @("E45", "  +2.52%  "),
    @("E46", "  -1.31%  "),
    @("E47", "  +6.74%  "),
    @("B48", "BabyDogeCoin"),
    @("C48", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D48", "0.0₆0351"),
    @("E48", "  -2.86%  "),
    @("B49", "Monero"),
    @("C49", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D49", "146.31"),
    @("E49", "  +1.80%  "),
    @("B50", "Maker"),
    @("C50", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"),
    @("D50", "2.751.96"),
    @("E50", "  +14.76%  "),
    @("B51", "FLOKI"),
    @("C51", "https://coinranking.com/coin/fmHk13Rqw+floki-floki"),
    @("D51", "0.000270"),
    @("E51", "  +7.95%  ")
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

